# edit.ps1 - applies the "Updated the initial planning meeting" change:
#   1. After "Please send it back soonest..." paragraph, add a blank
#      paragraph and a "Below are the high-level milestone dates." paragraph.
#   2. After "Warm regards, Pieter" paragraph, add a blank paragraph, an
#      (empty) paragraph carrying a minorHAnsi/minorBidi rFonts run-properties
#      mark, and a 2-column x 14-row milestone table.

$d = $word.ActiveDocument

function Insert-XmlAfterText([string]$searchText, [string]$xmlFragment) {
    $found = $d.Content
    $found.Find.ClearFormatting()
    $ok = $found.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Insert-XmlAfterText: text not found: $searchText"
    }
    $pos = $found.End
    $rng = $d.Range($pos, $pos)
    $rng.InsertXML($xmlFragment)
}

# --- Insertion 1: milestone intro sentence, right after the "scope approval" ask ---
$para1Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Below are the high-level milestone dates.</w:t></w:r></w:p>'
Insert-XmlAfterText "Please send it back soonest so that I can request scope approval." $para1Xml

# --- Insertion 2: blank paragraph + font marker paragraph + milestone table, after the sign-off ---
$blankParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$fontsParaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorBidi"/>
        </w:rPr>
      </w:pPr>
    </w:p>
'@
$tableXml = @'
<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:tblPr>
        <w:tblW w:w="5732" w:type="dxa"/>
        <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
      </w:tblPr>
      <w:tblGrid>
        <w:gridCol w:w="4316"/>
        <w:gridCol w:w="1416"/>
      </w:tblGrid>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="288"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="4316" w:type="dxa"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t>&lt;COMP Name&gt;</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t xml:space="preserve"> - Evaluation appraisal</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1416" w:type="dxa"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="right"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t>25-Sep-23</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="300"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="4316" w:type="dxa"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t>&lt;COMP Name&gt;</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t xml:space="preserve"> - Sponsor meeting</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1416" w:type="dxa"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="right"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t>02-Oct-23</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="300"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="4316" w:type="dxa"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t>&lt;COMP Name&gt;</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t xml:space="preserve"> - Scope approval</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1416" w:type="dxa"/>
            <w:tcBorders>
              <w:top w:val="single" w:sz="8" w:space="0" w:color="auto"/>
              <w:left w:val="single" w:sz="8" w:space="0" w:color="auto"/>
              <w:bottom w:val="single" w:sz="8" w:space="0" w:color="auto"/>
              <w:right w:val="single" w:sz="8" w:space="0" w:color="auto"/>
            </w:tcBorders>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="right"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t>09-Oct-23</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="288"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="4316" w:type="dxa"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t>&lt;COMP Name&gt;</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t xml:space="preserve"> - Scope must be done</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1416" w:type="dxa"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="right"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t>29-Oct-23</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="288"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="4316" w:type="dxa"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t>&lt;COMP Name&gt;</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t xml:space="preserve"> - 30% review</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1416" w:type="dxa"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="right"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t>10-Nov-23</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="288"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="4316" w:type="dxa"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t>&lt;COMP Name&gt;</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t xml:space="preserve"> - 50% review</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1416" w:type="dxa"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="right"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t>17-Nov-23</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="288"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="4316" w:type="dxa"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t>&lt;COMP Name&gt;</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t xml:space="preserve"> - 90% review</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1416" w:type="dxa"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="right"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t>28-Nov-23</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="288"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="4316" w:type="dxa"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t>&lt;COMP Name&gt;</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t xml:space="preserve"> - EPG meeting</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1416" w:type="dxa"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="right"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t>01-Dec-23</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="288"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="4316" w:type="dxa"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t>&lt;COMP Name&gt;</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t xml:space="preserve"> - Technology test</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1416" w:type="dxa"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="right"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t>06-Dec-23</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="288"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="4316" w:type="dxa"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t>&lt;COMP Name&gt;</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t xml:space="preserve"> - Phase 2 start</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1416" w:type="dxa"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="right"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t>08-Dec-23</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="288"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="4316" w:type="dxa"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t>&lt;COMP Name&gt;</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t xml:space="preserve"> - Phase 3 done</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1416" w:type="dxa"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="right"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t>15-Dec-23</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="288"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="4316" w:type="dxa"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t>&lt;COMP Name&gt;</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t xml:space="preserve"> - Finalise CAS</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1416" w:type="dxa"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="right"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t>18-Dec-23</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="288"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="4316" w:type="dxa"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t>&lt;COMP Name&gt;</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t xml:space="preserve"> - Final CAS review</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1416" w:type="dxa"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="right"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t>21-Dec-23</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="288"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="4316" w:type="dxa"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t>&lt;COMP Name&gt;</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:color w:val="000000"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t xml:space="preserve"> - Submit</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="1416" w:type="dxa"/>
            <w:noWrap/>
            <w:vAlign w:val="bottom"/>
            <w:hideMark/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:jc w:val="right"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman"/>
                <w:lang w:eastAsia="en-ZA"/>
              </w:rPr>
              <w:t>24-Dec-23</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
    </w:tbl>
'@

$secondFragment = $blankParaXml + $fontsParaXml + $tableXml
Insert-XmlAfterText "Warm regards, Pieter" $secondFragment

Write-Output "edit applied"
